# Resultado das instâncias da simulação.
# Adds a new "LCR" / "Antigo" mini-table (M1:N1 header) and a third
# algorithm column "GRASP" next to the existing "Guloso" / "Busca Local"
# columns (J/K), replicated into M/N/O, with matching data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Row 1: small "LCR" / "Antigo" header, styled like the J2/K2 header ---
$ws.Range("J2").Copy()
$ws.Range("M1:O1").PasteSpecial($xlPasteFormats)

$ws.Range("M1").Value = "LCR"
$ws.Range("N1").Value = "Antigo"

# --- Row 2: mirror the Guloso / Busca Local headers, plus new GRASP ---
$ws.Range("J2:K2").Copy()
$ws.Range("M2:N2").PasteSpecial($xlPasteFormats)
$ws.Range("J2").Copy()
$ws.Range("O2").PasteSpecial($xlPasteFormats)

$ws.Range("M2").Value = "Guloso"
$ws.Range("N2").Value = "Busca Local"
$ws.Range("O2").Value = "GRASP"

# --- Data rows 3-10: new M/N/O values, styled like the existing J/K data ---
$ws.Range("J3:K10").Copy()
$ws.Range("M3:N10").PasteSpecial($xlPasteFormats)
$ws.Range("J3").Copy()
$ws.Range("O3:O10").PasteSpecial($xlPasteFormats)

$ws.Range("M3").Value = 7
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 6

$ws.Range("M4").Value = 12
$ws.Range("N4").Value = 12
$ws.Range("O4").Value = 11

$ws.Range("M5").Value = 16
$ws.Range("N5").Value = 16
$ws.Range("O5").Value = 16

$ws.Range("M6").Value = 6
$ws.Range("N6").Value = 6
$ws.Range("O6").Value = 6

$ws.Range("M7").Value = 22
$ws.Range("N7").Value = 22
$ws.Range("O7").Value = 21

$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 20
$ws.Range("O8").Value = 20

$ws.Range("M9").Value = 29
$ws.Range("N9").Value = 28
$ws.Range("O9").Value = 28

$ws.Range("M10").Value = 21
$ws.Range("N10").Value = 21
$ws.Range("O10").Value = 21

# Touch a block of otherwise-empty trailing rows (no cell content), matching
# the incidental row bookkeeping left behind in the authored workbook.
for ($r = 13; $r -le 28; $r++) {
    $ws.Rows($r).RowHeight = 13.8
}

# Restore the selection to roughly where it ended up after the edit.
$ws.Range("N15").Select()
